$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) figures per latest data refresh

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.160.52'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.628.02'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.585'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.627.57'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.108'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.60'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.362'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.16'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.098.16'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.007.98'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.558.12'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.29'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.49'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '339.80'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.89'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.64%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  -4.20%  '
$ws.Range('E25').Value = '  -2.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.52'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.163'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.65%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '540.06'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.84'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.62%  '
$ws.Range('E33').Value = '  +1.23%  '
$ws.Range('E34').Value = '  -3.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0801'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.34'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +14.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '167.63'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.402'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.98'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.88'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.22%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '168.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.74'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.28'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0565'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.623'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0241'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0957'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.55'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.75'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.16%  '
